$wb = $excel.ActiveWorkbook

$ws_ALC = $wb.Worksheets.Item("ALC")
$ws_ALC.Range("H19").Value = 1590.0834
$ws_ALC.Range("I19").Value = 801
$ws_ALC.Range("J19").Value = 1661.8182
$ws_ALC.Range("K19").Value = 801
$ws_ALC.Range("L19").Value = 1661.8182
$ws_ALC.Range("M19").Value = -626
$ws_ALC.Range("N19").Value = -2011.8182
$ws_ALC.Range("H28").Value = 553.08826
$ws_ALC.Range("I28").Value = 465.45834
$ws_ALC.Range("J28").Value = 763.4
$ws_ALC.Range("K28").Value = 465.45834
$ws_ALC.Range("L28").Value = 763.4
$ws_ALC.Range("M28").Value = 19.54165999999998
$ws_ALC.Range("N28").Value = -1733.4
$ws_ALC.Range("H40").Value = 1497
$ws_ALC.Range("I40").Value = 1000
$ws_ALC.Range("J40").Value = 1994
$ws_ALC.Range("K40").Value = 1000
$ws_ALC.Range("L40").Value = 1994
$ws_ALC.Range("M40").Value = -825
$ws_ALC.Range("N40").Value = -2344
$ws_ALC.Range("H43").Value = 5680.615
$ws_ALC.Range("J43").Value = 4651.25
$ws_ALC.Range("L43").Value = 4651.25
$ws_ALC.Range("N43").Value = -4789.25
$ws_ALC.Range("H112").Value = 2008.7273
$ws_ALC.Range("J112").Value = 1847.3334
$ws_ALC.Range("L112").Value = 5542.0002
$ws_ALC.Range("N112").Value = -7758.0002
$ws_ALC.Range("H137").Value = 10534836
$ws_ALC.Range("I137").Value = 16668692
$ws_ALC.Range("J137").Value = 19656.285
$ws_ALC.Range("K137").Value = 50006076
$ws_ALC.Range("L137").Value = 58968.855
$ws_ALC.Range("M137").Value = -50003526
$ws_ALC.Range("N137").Value = -64068.855
$ws_ALC.Range("H138").Value = 6261.0757
$ws_ALC.Range("I138").Value = 4977.684
$ws_ALC.Range("J138").Value = 6978.2646
$ws_ALC.Range("K138").Value = 14933.052
$ws_ALC.Range("L138").Value = 20934.7938
$ws_ALC.Range("M138").Value = -9793.052
$ws_ALC.Range("N138").Value = -31214.7938
$ws_ARM = $wb.Worksheets.Item("ARM")
$ws_ARM.Range("H32").Value = 784754.0600000001
$ws_ARM.Range("I32").Value = 795979.5600000001
$ws_ARM.Range("K32").Value = 795979.5600000001
$ws_ARM.Range("M32").Value = -795692.5600000001
$ws_ARM.Range("H61").Value = 8302009.5
$ws_ARM.Range("I61").Value = 4133918.8
$ws_ARM.Range("K61").Value = 4133918.8
$ws_ARM.Range("M61").Value = -4133706.8
$ws_ARM.Range("H74").Value = 932537.75
$ws_ARM.Range("I74").Value = 1246606
$ws_ARM.Range("K74").Value = 1246606
$ws_ARM.Range("M74").Value = -1245732
$ws_ARM.Range("H77").Value = 932537.75
$ws_ARM.Range("I77").Value = 1246606
$ws_ARM.Range("K77").Value = 6233030
$ws_ARM.Range("M77").Value = -6228662
$ws_ARM.Range("H97").Value = 865.1923
$ws_ARM.Range("I97").Value = 192.2381
$ws_ARM.Range("J97").Value = 3691.6
$ws_ARM.Range("K97").Value = 192.2381
$ws_ARM.Range("L97").Value = 3691.6
$ws_ARM.Range("M97").Value = 303.7619
$ws_ARM.Range("N97").Value = -4683.6
$ws_ARM.Range("H125").Value = 0
$ws_ARM.Range("I125").Value = 0
$ws_ARM.Range("K125").Value = 0
$ws_ARM.Range("M125").ClearContents()
$ws_ARM.Range("H132").Value = 3175.7192
$ws_ARM.Range("I132").Value = 1846.8684
$ws_ARM.Range("K132").Value = 5540.6052
$ws_ARM.Range("M132").Value = -3010.6052
$ws_ARM.Range("H136").Value = 8302009.5
$ws_ARM.Range("I136").Value = 4133918.8
$ws_ARM.Range("K136").Value = 12401756.4
$ws_ARM.Range("M136").Value = -12399206.4
$ws_BSM = $wb.Worksheets.Item("BSM")
$ws_BSM.Range("H20").Value = 761264.3
$ws_BSM.Range("I20").Value = 909347.1
$ws_BSM.Range("K20").Value = 909347.1
$ws_BSM.Range("M20").Value = -909100.1
$ws_BSM.Range("H26").Value = 32999.75
$ws_BSM.Range("I26").Value = 22333
$ws_BSM.Range("K26").Value = 22333
$ws_BSM.Range("M26").Value = -22041
$ws_BSM.Range("H134").Value = 8638711
$ws_BSM.Range("I134").Value = 7960301.5
$ws_BSM.Range("J134").Value = 13896386
$ws_BSM.Range("K134").Value = 23880904.5
$ws_BSM.Range("L134").Value = 41689158
$ws_BSM.Range("M134").Value = -23878369.5
$ws_BSM.Range("N134").Value = -41694228
$ws_CRP = $wb.Worksheets.Item("CRP")
$ws_CRP.Range("H22").Value = 816.4815
$ws_CRP.Range("I22").Value = 602.8261
$ws_CRP.Range("K22").Value = 602.8261
$ws_CRP.Range("M22").Value = -252.8261
$ws_CRP.Range("H31").Value = 563379.5
$ws_CRP.Range("I31").Value = 851637.8
$ws_CRP.Range("K31").Value = 851637.8
$ws_CRP.Range("M31").Value = -851342.8
$ws_CRP.Range("H34").Value = 563379.5
$ws_CRP.Range("I34").Value = 851637.8
$ws_CRP.Range("K34").Value = 851637.8
$ws_CRP.Range("M34").Value = -851435.8
$ws_CRP.Range("H58").Value = 7008451.5
$ws_CRP.Range("I58").Value = 13892330
$ws_CRP.Range("K58").Value = 13892330
$ws_CRP.Range("M58").Value = -13892127
$ws_CRP.Range("H132").Value = 1559.5
$ws_CRP.Range("I132").Value = 1274.375
$ws_CRP.Range("J132").Value = 2700
$ws_CRP.Range("K132").Value = 3823.125
$ws_CRP.Range("L132").Value = 8100
$ws_CRP.Range("M132").Value = -1293.125
$ws_CRP.Range("N132").Value = -13160
$ws_CRP.Range("H134").Value = 4009.377
$ws_CRP.Range("I134").Value = 2325.1785
$ws_CRP.Range("K134").Value = 6975.5355
$ws_CRP.Range("M134").Value = -4440.5355
$ws_CRP.Range("H136").Value = 7008451.5
$ws_CRP.Range("I136").Value = 13892330
$ws_CRP.Range("K136").Value = 41676990
$ws_CRP.Range("M136").Value = -41674440
$ws_CUL = $wb.Worksheets.Item("CUL")
$ws_CUL.Range("H58").Value = 12499.875
$ws_CUL.Range("I58").Value = 9999
$ws_CUL.Range("K58").Value = 29997
$ws_CUL.Range("M58").Value = -29869
$ws_CUL.Range("H60").Value = 5043.5713
$ws_CUL.Range("I60").Value = 5866.8335
$ws_CUL.Range("K60").Value = 17600.5005
$ws_CUL.Range("M60").Value = -17349.5005
$ws_CUL.Range("H122").Value = 1153007.6
$ws_CUL.Range("I122").Value = 2016828.4
$ws_CUL.Range("J122").Value = 1246.6666
$ws_CUL.Range("K122").Value = 18151455.6
$ws_CUL.Range("L122").Value = 11219.9994
$ws_CUL.Range("M122").Value = -18149005.6
$ws_CUL.Range("N122").Value = -16119.9994
$ws_GSM = $wb.Worksheets.Item("GSM")
$ws_GSM.Range("H2").Value = 4590958.5
$ws_GSM.Range("I2").Value = 6312541
$ws_GSM.Range("K2").Value = 6312541
$ws_GSM.Range("M2").Value = -6312428
$ws_GSM.Range("H80").Value = 3289.6365
$ws_GSM.Range("I80").Value = 3242.353
$ws_GSM.Range("J80").Value = 3450.4
$ws_GSM.Range("K80").Value = 3242.353
$ws_GSM.Range("L80").Value = 3450.4
$ws_GSM.Range("M80").Value = -2244.353
$ws_GSM.Range("N80").Value = -5446.4
$ws_GSM.Range("H83").Value = 3289.6365
$ws_GSM.Range("I83").Value = 3242.353
$ws_GSM.Range("J83").Value = 3450.4
$ws_GSM.Range("K83").Value = 16211.765
$ws_GSM.Range("L83").Value = 17252
$ws_GSM.Range("M83").Value = -11219.765
$ws_GSM.Range("N83").Value = -27236
$ws_GSM.Range("H97").Value = 754.1539
$ws_GSM.Range("I97").Value = 668.9268
$ws_GSM.Range("J97").Value = 1071.8182
$ws_GSM.Range("K97").Value = 668.9268
$ws_GSM.Range("L97").Value = 1071.8182
$ws_GSM.Range("M97").Value = -172.9268
$ws_GSM.Range("N97").Value = -2063.8182
$ws_GSM.Range("H102").Value = 2426.3257
$ws_GSM.Range("I102").Value = 1984.5
$ws_GSM.Range("J102").Value = 3711.6365
$ws_GSM.Range("K102").Value = 1984.5
$ws_GSM.Range("L102").Value = 3711.6365
$ws_GSM.Range("M102").Value = -362.5
$ws_GSM.Range("N102").Value = -6955.636500000001
$ws_GSM.Range("H113").Value = 1421.1428
$ws_GSM.Range("I113").Value = 1490.2222
$ws_GSM.Range("J113").Value = 1006.6667
$ws_GSM.Range("K113").Value = 1490.2222
$ws_GSM.Range("L113").Value = 1006.6667
$ws_GSM.Range("M113").Value = 679.7778000000001
$ws_GSM.Range("N113").Value = -5346.6667
$ws_GSM.Range("H132").Value = 14516.8
$ws_GSM.Range("I132").Value = 10175.762
$ws_GSM.Range("K132").Value = 30527.286
$ws_GSM.Range("M132").Value = -27997.286
$ws_LTW = $wb.Worksheets.Item("LTW")
$ws_LTW.Range("H22").Value = 2851.8572
$ws_LTW.Range("I22").Value = 2216.3333
$ws_LTW.Range("J22").Value = 3106.0667
$ws_LTW.Range("K22").Value = 2216.3333
$ws_LTW.Range("L22").Value = 3106.0667
$ws_LTW.Range("M22").Value = -1921.3333
$ws_LTW.Range("N22").Value = -3696.0667
$ws_LTW.Range("H27").Value = 2851.8572
$ws_LTW.Range("I27").Value = 2216.3333
$ws_LTW.Range("J27").Value = 3106.0667
$ws_LTW.Range("K27").Value = 2216.3333
$ws_LTW.Range("L27").Value = 3106.0667
$ws_LTW.Range("M27").Value = -2109.3333
$ws_LTW.Range("N27").Value = -3320.0667
$ws_LTW.Range("H46").Value = 4308.4
$ws_LTW.Range("I46").Value = 2320.25
$ws_LTW.Range("J46").Value = 6143.615
$ws_LTW.Range("K46").Value = 2320.25
$ws_LTW.Range("L46").Value = 6143.615
$ws_LTW.Range("M46").Value = -2132.25
$ws_LTW.Range("N46").Value = -6519.615
$ws_LTW.Range("H61").Value = 6830.5757
$ws_LTW.Range("I61").Value = 6181
$ws_LTW.Range("K61").Value = 6181
$ws_LTW.Range("M61").Value = -5979
$ws_LTW.Range("H100").Value = 2552
$ws_LTW.Range("I100").Value = 2260.6667
$ws_LTW.Range("K100").Value = 2260.6667
$ws_LTW.Range("M100").Value = -1719.6667
$ws_LTW.Range("H113").Value = 6830.5757
$ws_LTW.Range("I113").Value = 6181
$ws_LTW.Range("K113").Value = 6181
$ws_LTW.Range("M113").Value = -4011
$ws_LTW.Range("H132").Value = 4633096.5
$ws_LTW.Range("I132").Value = 6413429.5
$ws_LTW.Range("K132").Value = 19240288.5
$ws_LTW.Range("M132").Value = -19237758.5
$ws_WVR = $wb.Worksheets.Item("WVR")
$ws_WVR.Range("H100").Value = 608.8
$ws_WVR.Range("I100").Value = 680.8333
$ws_WVR.Range("J100").Value = 500.75
$ws_WVR.Range("K100").Value = 1361.6666
$ws_WVR.Range("L100").Value = 1001.5
$ws_WVR.Range("M100").Value = -820.6666
$ws_WVR.Range("N100").Value = -2083.5
$ws_WVR.Range("H113").Value = 1080.3611
$ws_WVR.Range("I113").Value = 699.84
$ws_WVR.Range("J113").Value = 1945.1818
$ws_WVR.Range("K113").Value = 2099.52
$ws_WVR.Range("L113").Value = 5835.5454
$ws_WVR.Range("M113").Value = 70.48000000000002
$ws_WVR.Range("N113").Value = -10175.5454
$ws_WVR.Range("H122").Value = 46713.72
$ws_WVR.Range("I122").Value = 1707
$ws_WVR.Range("J122").Value = 189235
$ws_WVR.Range("K122").Value = 5121
$ws_WVR.Range("L122").Value = 567705
$ws_WVR.Range("M122").Value = -2671
$ws_WVR.Range("N122").Value = -572605
$ws_WVR.Range("H125").Value = 59969
$ws_WVR.Range("J125").Value = 59969
$ws_WVR.Range("L125").Value = 59969
$ws_WVR.Range("N125").Value = -69809
$ws_WVR.Range("H132").Value = 5748687.5
$ws_WVR.Range("I132").Value = 8772869
$ws_WVR.Range("J132").Value = 2742.9
$ws_WVR.Range("K132").Value = 26318607
$ws_WVR.Range("L132").Value = 8228.700000000001
$ws_WVR.Range("M132").Value = -26316077
$ws_WVR.Range("N132").Value = -13288.7
$ws_WVR.Range("H136").Value = 3534373
$ws_WVR.Range("I136").Value = 2899846.2
$ws_WVR.Range("K136").Value = 8699538.600000001
$ws_WVR.Range("M136").Value = -8699538.600000001
